$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new "2021年" data row (row 11) below the existing data (which
# currently ends at row 10, "2020年"), matching the columns laid out in row 1.

$ws.Range("A11").Value = "2021年"
$ws.Range("B11").Value = 2138
$ws.Range("C11").Value = 724
$ws.Range("D11").Value = 188
# Column E has no figure for this year (same as E10) - write it as an empty
# text cell rather than leaving it completely blank.
$ws.Range("E11").Value = "'"
$ws.Range("E11").Style = "Normal"
$ws.Range("F11").Value = 1811
$ws.Range("G11").Value = 2315
$ws.Range("H11").Value = 271
$ws.Range("I11").Value = 1603
$ws.Range("J11").Value = 560
$ws.Range("K11").Value = 47427
$ws.Range("L11").Value = 711
$ws.Range("M11").Value = 99
$ws.Range("N11").Value = 58
$ws.Range("O11").Value = 1143
$ws.Range("P11").Value = 976
$ws.Range("Q11").Value = 231
$ws.Range("R11").Value = 335
$ws.Range("S11").Value = 1662
$ws.Range("T11").Value = 393
$ws.Range("U11").Value = 2877
$ws.Range("V11").Value = 71
$ws.Range("W11").Value = 1414
$ws.Range("X11").Value = 213
$ws.Range("Y11").Value = 1153
$ws.Range("Z11").Value = 3791
$ws.Range("AA11").Value = 1063
$ws.Range("AB11").Value = 497
$ws.Range("AC11").Value = 49
$ws.Range("AD11").Value = 1905
$ws.Range("AE11").Value = 1775
$ws.Range("AF11").Value = 5081
$ws.Range("AG11").Value = 2562
$ws.Range("AH11").Value = 640
$ws.Range("AI11").Value = 664
$ws.Range("AJ11").Value = 145
$ws.Range("AK11").Value = 2188
$ws.Range("AL11").Value = 879
$ws.Range("AM11").Value = 2771
$ws.Range("AN11").Value = 111
$ws.Range("AO11").Value = 1318
$ws.Range("AP11").Value = 849
$ws.Range("AQ11").Value = 192

# Match the row-label formatting (bold, bordered, centered) used by the
# other year cells in column A.
$ws.Range("A10").Copy()
$ws.Range("A11").PasteSpecial(-4122)
